# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
#
# Refresh the "MSME Country Indicators - Luxembourg Summary" sheet with
# updated two-decimal-precision figures for the SME Associations source
# block (Enterprises density, Employment %, Enterprises %, Value added %).
# Values are stored as text (as in the source data), so a leading
# apostrophe is used to force text entry without altering number formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Source Type: Statistical Institution - Enterprises density (per 1000 people)
$ws.Range("B11").Value = "'54.64"
$ws.Range("C11").Value = "'8.18"
$ws.Range("D11").Value = "'62.82"

# Source Type: SME Associations - Enterprises density (per 1000 people)
$ws.Range("B33").Value = "'48.97"
$ws.Range("C33").Value = "'7.02"
$ws.Range("D33").Value = "'55.99"

# Source Type: SME Associations - Employment (% of total)
$ws.Range("B34").Value = "'23.11"
$ws.Range("C34").Value = "'42.73"
$ws.Range("D34").Value = "'65.84"

# Source Type: SME Associations - Enterprises (% of total)
$ws.Range("B36").Value = "'87.03"
$ws.Range("C36").Value = "'12.48"
$ws.Range("D36").Value = "'99.51"

# Source Type: SME Associations - Value added to the economy (% of total)
$ws.Range("B40").Value = "'29.88"
$ws.Range("C40").Value = "'39.28"
$ws.Range("D40").Value = "'69.17"
